# Auto-generated Excel COM-interop edit script
# Applies the "cryptos" price-list refresh described by the commit:
#   - Price (D) and Volume 1h (E) values are updated for most rows.
#   - A handful of rows swap places (their rank changed), so for
#     those the whole row (Coin/Link/Price/Volume) is rewritten.
# D holds values such as "307.49" or "44.174.30" as plain text in the
# source workbook (inline strings). Excel's automatic type detection
# would otherwise coerce plain numeric-looking text to a Number, so we
# force the Text number format on any D cell we touch before writing,
# keeping the value exactly as authored.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.174.30"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.245.75"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.49"
$ws.Range("E5").Value = "  -1.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.22"
$ws.Range("E6").Value = "  -2.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  +1.07%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -1.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.00"
$ws.Range("E10").Value = "  -3.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.28"
$ws.Range("E12").Value = "  -1.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.584.69"
$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.243.55"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.834"
$ws.Range("E16").Value = "  -0.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.66"
$ws.Range("E17").Value = "  -3.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.050.69"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0974"
$ws.Range("E19").Value = "  +1.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.29"
$ws.Range("E20").Value = "  -5.86%  "

$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.46"
$ws.Range("E22").Value = "  +0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.98"
$ws.Range("E23").Value = "  +1.50%  "

$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").Value = "  -1.75%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.95"
$ws.Range("E27").Value = "  -2.59%  "

$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.31"
$ws.Range("E28").Value = "  +3.93%  "

$ws.Range("E29").Value = "  +1.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.05"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.10"
$ws.Range("E31").Value = "  +0.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.33"
$ws.Range("E32").Value = "  -3.59%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0803"
$ws.Range("E33").Value = "  -3.28%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.32"
$ws.Range("E34").Value = "  +4.90%  "

$ws.Range("E35").Value = "  -3.52%  "

$ws.Range("E36").Value = "  +2.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -6.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.87"
$ws.Range("E39").Value = "  -6.64%  "

$ws.Range("E40").Value = "  -5.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.87"
$ws.Range("E41").Value = "  -4.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0298"
$ws.Range("E42").Value = "  -3.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.723.64"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "83.77"
$ws.Range("E45").Value = "  +3.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.191"
$ws.Range("E46").Value = "  -1.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.46"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.16"
$ws.Range("E48").Value = "  +1.29%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.86"
$ws.Range("E49").Value = "  -4.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.03"
$ws.Range("E50").Value = "  -5.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.70"
$ws.Range("E51").Value = "  -3.05%  "
